# Roboflow Annotation Report 7/3/2025
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in row 50 with this week's annotation progress data
$ws.Range("D50").Value = 45723
$ws.Range("E50").Value = 192
$ws.Range("F50").Value = 734
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 1012
$ws.Range("J50").Value = "N/A"

# Move the active selection as left after editing
[void]$ws.Range("J55").Select()
